$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.452.14'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '3.561.70'
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.25'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.94'
$ws.Range("E6").Value = '  +0.79%  '

$ws.Range("D7").Value = '3.560.60'
$ws.Range("E7").Value = '  +0.72%  '

$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("E9").Value = '  +3.75%  '

$ws.Range("E10").Value = '  -0.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.93'
$ws.Range("E11").Value = '  -1.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.414'
$ws.Range("E12").Value = '  +0.78%  '

$ws.Range("D13").Value = '4.159.72'
$ws.Range("E13").Value = '  +0.58%  '

$ws.Range("E14").Value = '  -0.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.03'
$ws.Range("E15").Value = '  -0.71%  '

$ws.Range("D16").Value = '3.561.98'
$ws.Range("E16").Value = '  +0.84%  '

$ws.Range("D17").Value = '66.509.97'
$ws.Range("E17").Value = '  +0.21%  '

$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.60'
$ws.Range("E19").Value = '  +6.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.18'
$ws.Range("E20").Value = '  -0.47%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.82'
$ws.Range("E21").Value = '  -0.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '431.56'
$ws.Range("E22").Value = '  +1.35%  '

$ws.Range("E23").Value = '  +1.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.76'
$ws.Range("E24").Value = '  +1.40%  '

$ws.Range("D25").Value = '3.703.44'
$ws.Range("E25").Value = '  +0.76%  '

$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("E27").Value = '  -0.75%  '

$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.51'
$ws.Range("E28").Value = '  +1.16%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.00'
$ws.Range("E29").Value = '  -1.28%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.18'
$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.45'
$ws.Range("E32").Value = '  -2.32%  '

$ws.Range("D33").Value = '3.554.37'
$ws.Range("E33").Value = '  +0.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.38'
$ws.Range("E34").Value = '  +0.44%  '

$ws.Range("E35").Value = '  -4.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.87'
$ws.Range("E36").Value = '  +0.77%  '

$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("E38").Value = '  -2.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.59'
$ws.Range("E39").Value = '  -0.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '174.74'
$ws.Range("E40").Value = '  +1.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0849'
$ws.Range("E41").Value = '  -1.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.20'
$ws.Range("E42").Value = '  +0.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.889'
$ws.Range("E43").Value = '  -0.29%  '

$ws.Range("E44").Value = '  +1.33%  '

$ws.Range("E45").Value = '  +1.28%  '

$ws.Range("E46").Value = '  +0.03%  '

$ws.Range("E47").Value = '  +4.56%  '

$ws.Range("E48").Value = '  -2.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.12'
$ws.Range("E49").Value = '  -3.91%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.71'
$ws.Range("E50").Value = '  +5.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.15'
$ws.Range("E51").Value = '  +0.14%  '
